# Applies the crypto price/volume refresh described in the commit diff.
# Columns B/C (text) and E (percent strings) are plain text naturally;
# column D ("Price") contains dotted-thousands strings like "41.511.34" that
# Excel would otherwise try to coerce into numbers/dates, so we force those
# cells to Text format before assigning, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCells = @("D2", "D3", "D5", "D6", "D10", "D11", "D13", "D15", "D16", "D17", "D18", "D19", "D21", "D23", "D24", "D27", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D39", "D41", "D43", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.511.34'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '2.488.69'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '315.30'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '94.29'
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('E9').Value = '  -3.09%  '
$ws.Range('D10').Value = '33.60'
$ws.Range('E10').Value = '  -4.33%  '
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '2.872.18'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('E14').Value = '  -3.43%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.50'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.493.89'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '41.469.73'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').Value = '6.35'
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('D21').Value = '70.19'
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('E22').Value = '  -6.53%  '
$ws.Range('D23').Value = '236.47'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').Value = '2.78'
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  -4.27%  '
$ws.Range('D27').Value = '24.20'
$ws.Range('E27').Value = '  -4.26%  '
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = '37.11'
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('D31').Value = '154.49'
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').Value = '  -5.23%  '
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('D34').Value = '0.0756'
$ws.Range('E34').Value = '  -3.05%  '
$ws.Range('D35').Value = '17.95'
$ws.Range('E35').Value = '  +2.98%  '
$ws.Range('D36').Value = '3.07'
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').Value = '2.42'
$ws.Range('E37').Value = '  -10.66%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('D39').Value = '0.114'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('E40').Value = '  -5.66%  '
$ws.Range('D41').Value = '4.14'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').Value = '19.84'
$ws.Range('E43').Value = '  -7.61%  '
$ws.Range('D44').Value = '1.991.50'
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '0.0286'
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('E46').Value = '  -5.61%  '
$ws.Range('D47').Value = '8.84'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '2.734.78'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').Value = '69.56'
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('D50').Value = '97.26'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('E51').Value = '  -4.70%  '
